# Update the K column (G) with recalculated strikeout counts
# (regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newK = @{
    2  = 4
    3  = 3
    4  = 3
    5  = 4
    6  = 3
    7  = 6
    8  = 4
    9  = 2
    10 = 16
    11 = 6
    12 = 1
    13 = 6
    14 = 10
    15 = 8
    16 = 5
    17 = 2
    18 = 5
    19 = 6
    20 = 4
    21 = 6
    22 = 5
    23 = 2
    24 = 4
    25 = 3
    26 = 1
    27 = 3
}

foreach ($row in $newK.Keys) {
    $ws.Range("G$row").Value = $newK[$row]
}
